$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from row 709 (last existing styled row) so the
# newly appended rows 710-727 reuse the same cell styles (s=23 / s=24)
# instead of Excel fabricating brand-new font/fill entries.
$fmtSource = $ws.Range("A709:C709")

$fmtSource.Copy()
$ws.Range("A710:C710").PasteSpecial(-4122)
$ws.Range("A710").Value = 'statistics.femaleChild'
$ws.Range("B710").Value = 'Female Child'
$ws.Range("C710").Value = '&Female Child'

$fmtSource.Copy()
$ws.Range("A711:C711").PasteSpecial(-4122)
$ws.Range("A711").Value = 'statistics.maleChild'
$ws.Range("B711").Value = 'Male Child'
$ws.Range("C711").Value = '&Male Child'

$fmtSource.Copy()
$ws.Range("A712:C712").PasteSpecial(-4122)
$ws.Range("A712").Value = 'statistics.femaleAdult'
$ws.Range("B712").Value = 'Female Adult'
$ws.Range("C712").Value = '&Female Adult'

$fmtSource.Copy()
$ws.Range("A713:C713").PasteSpecial(-4122)
$ws.Range("A713").Value = 'statistics.maleAdult'
$ws.Range("B713").Value = 'Male Adult'
$ws.Range("C713").Value = '&Male Adult'

$fmtSource.Copy()
$ws.Range("A714:C714").PasteSpecial(-4122)
$ws.Range("A714").Value = 'statistics.totalFChild'
$ws.Range("B714").Value = 'Total Female Children:'
$ws.Range("C714").Value = '&Total Female Children:'

$fmtSource.Copy()
$ws.Range("A715:C715").PasteSpecial(-4122)
$ws.Range("A715").Value = 'statistics.totalMChild'
$ws.Range("B715").Value = 'Total Male Children:'
$ws.Range("C715").Value = '&Total Male Children:'

$fmtSource.Copy()
$ws.Range("A716:C716").PasteSpecial(-4122)
$ws.Range("A716").Value = 'statistics.totalFAdult'
$ws.Range("B716").Value = 'Total Female Adults:'
$ws.Range("C716").Value = '&Total Female Adults:'

$fmtSource.Copy()
$ws.Range("A717:C717").PasteSpecial(-4122)
$ws.Range("A717").Value = 'statistics.totalMAdult'
$ws.Range("B717").Value = 'Total Male Adults:'
$ws.Range("C717").Value = '&Total Male Adults:'

$fmtSource.Copy()
$ws.Range("A718:C718").PasteSpecial(-4122)
$ws.Range("A718").Value = 'statistics.totalFChildFollowUpVisits'
$ws.Range("B718").Value = 'Total Female Children Follow Up Visits:'
$ws.Range("C718").Value = '&Total Female Children Follow Up Visits:'

$fmtSource.Copy()
$ws.Range("A719:C719").PasteSpecial(-4122)
$ws.Range("A719").Value = 'statistics.totalMChildFollowUpVisits'
$ws.Range("B719").Value = 'Total Male Children Follow Up Visits:'
$ws.Range("C719").Value = '&Total Male Children Follow Up Visits:'

$fmtSource.Copy()
$ws.Range("A720:C720").PasteSpecial(-4122)
$ws.Range("A720").Value = 'statistics.totalFAdultFollowUpVisits'
$ws.Range("B720").Value = 'Total Female Adult Follow Up Visits:'
$ws.Range("C720").Value = '&Total Female Adult Follow Up Visits:'

$fmtSource.Copy()
$ws.Range("A721:C721").PasteSpecial(-4122)
$ws.Range("A721").Value = 'statistics.totalMAdultFollowUpVisits'
$ws.Range("B721").Value = 'Total Male Adult Follow Up Visits:'
$ws.Range("C721").Value = '&Total Male Adult Follow Up Visits:'

$fmtSource.Copy()
$ws.Range("A722:C722").PasteSpecial(-4122)
$ws.Range("A722").Value = 'statistics.selectAtLeastOne'
$ws.Range("B722").Value = 'Select at least one Gender and Age option '
$ws.Range("C722").Value = '&Select at least one Gender and Age option '

$fmtSource.Copy()
$ws.Range("A723:C723").PasteSpecial(-4122)
$ws.Range("A723").Value = 'statistics.warning'
$ws.Range("B723").Value = 'Warning'
$ws.Range("C723").Value = '&Warning'

$fmtSource.Copy()
$ws.Range("A724:C724").PasteSpecial(-4122)
$ws.Range("A724").Value = 'statistics.totalDisFChild'
$ws.Range("B724").Value = 'Total Female Children With Disabilities: '
$ws.Range("C724").Value = '&Total Female Children With Disabilities: '

$fmtSource.Copy()
$ws.Range("A725:C725").PasteSpecial(-4122)
$ws.Range("A725").Value = 'statistics.totalDisMChild'
$ws.Range("B725").Value = 'Total Male Chidlren With Disabilities: '
$ws.Range("C725").Value = '&Total Male Chidlren With Disabilities: '

$fmtSource.Copy()
$ws.Range("A726:C726").PasteSpecial(-4122)
$ws.Range("A726").Value = 'statistics.totalDisFAdult'
$ws.Range("B726").Value = 'Total Female Adults With Disabilities: '
$ws.Range("C726").Value = '&Total Female Adults With Disabilities: '

$fmtSource.Copy()
$ws.Range("A727:C727").PasteSpecial(-4122)
$ws.Range("A727").Value = 'statistics.totalDisMAdult'
$ws.Range("B727").Value = 'Total Male Adults With Disabilities: '
$ws.Range("C727").Value = '&Total Male Adults With Disabilities: '
